# Append 7 new DTU "GetConfigReq" log rows (21-27) to the sheet, mirroring the
# layout of the existing rows. Only the timestamp columns (A - human readable
# datetime, AT - unix epoch seconds) and the wifi RSSI column (BF) change
# between rows; every other column repeats the same device/config snapshot
# that is already present on every prior row (passive/default CMNET config).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (1-based) that are identical across every row, with their common
# value. Column A (datetime), AT (unix time) and BF (wifi RSSI) are supplied
# per-row below and are intentionally excluded here.
$commonCols = [ordered]@{
    2  = ""                      # unknown
    3  = 0                       # accessModel
    4  = "NONE"                  # apnName
    5  = "NONE"                  # apnPassword
    6  = "CMNET"                 # apnSet
    7  = 223                     # cableDns0
    8  = 5                       # cableDns1
    9  = 5                       # cableDns2
    10 = 5                       # cableDns3
    11 = 0                       # channelSelect
    12 = 0                       # defaultGateway0
    13 = 0                       # defaultGateway1
    14 = 0                       # defaultGateway2
    15 = 0                       # defaultGateway3
    16 = 1                       # dhcpSwitch
    17 = "10F872226797"          # dtuSn
    18 = 0                       # invType
    19 = 0                       # ipAddr0
    20 = 0                       # ipAddr1
    21 = 0                       # ipAddr2
    22 = 0                       # ipAddr3
    23 = ""                      # kaNub
    24 = 0                       # limitPowerMyPower
    25 = 0                       # lockPassword
    26 = 0                       # lockTime
    27 = 72                      # mac0
    28 = 77                      # mac1
    29 = 114                     # mac2
    30 = 34                      # mac3
    31 = 103                     # mac4
    32 = 151                     # mac5
    33 = "NONE"                  # meterInterface
    34 = "NONE"                  # meterKind
    35 = 1                       # netmodeSelect
    36 = 3600                    # offset
    37 = "dataeu.hoymiles.com"   # serverDomainName
    38 = 1                       # serverSendTime
    39 = 10081                   # serverPort
    40 = 0                       # sub1GSweepSwitch
    41 = 0                       # sub1GWorkChannel
    42 = 0                       # subnetMask0
    43 = 0                       # subnetMask1
    44 = 0                       # subnetMask2
    45 = 0                       # subnetMask3
    47 = 0                       # wifiIpAddr0
    48 = 0                       # wifiIpAddr1
    49 = 0                       # wifiIpAddr2
    50 = 0                       # wifiIpAddr3
    51 = 0                       # wifiMac0
    52 = 0                       # wifiMac1
    53 = 0                       # wifiMac2
    54 = 0                       # wifiMac3
    55 = 0                       # wifiMac4
    56 = 0                       # wifiMac5
    57 = "0negawsklov0negawsklov" # wifiPassword
    59 = "HomeSweetHome"         # wifiSsid
    60 = 0                       # zeroExport433Addr
    61 = 0                       # zeroExportEnable
}

# Per-row data: A = datetime text, AT = unix epoch seconds, BF = wifi RSSI.
$newRows = @(
    @{ Row = 21; Dt = "2022-03-21 09:38:58"; Ts = 1647851939; Rssi = 58 },
    @{ Row = 22; Dt = "2022-03-21 09:58:00"; Ts = 1647853081; Rssi = 58 },
    @{ Row = 23; Dt = "2022-03-21 09:59:23"; Ts = 1647853164; Rssi = 62 },
    @{ Row = 24; Dt = "2022-03-21 10:01:22"; Ts = 1647853284; Rssi = 60 },
    @{ Row = 25; Dt = "2022-03-21 10:02:12"; Ts = 1647853334; Rssi = 66 },
    @{ Row = 26; Dt = "2022-03-21 10:03:18"; Ts = 1647853400; Rssi = 62 },
    @{ Row = 27; Dt = "2022-03-21 10:55:47"; Ts = 1647856549; Rssi = 68 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.Dt

    foreach ($col in $commonCols.Keys) {
        $ws.Cells.Item($r, $col).Value = $commonCols[$col]
    }

    $ws.Cells.Item($r, 46).Value = $entry.Ts   # time
    $ws.Cells.Item($r, 58).Value = $entry.Rssi # wifiRssi
}
